$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume/1h) store plain text values (e.g. "582.83",
# "  +2.08%  ") in the source data, not numbers/percentages. Force the Text
# number format first so Excel does not auto-convert the assigned strings into
# numeric values (which would also introduce floating point drift).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '68.071.31'
$ws.Range('E2').Value = '  +2.35%  '
$ws.Range('D3').Value = '3.343.55'
$ws.Range('E3').Value = '  +2.60%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '583.94'
$ws.Range('E5').Value = '  +3.15%  '
$ws.Range('D6').Value = '177.51'
$ws.Range('E6').Value = '  +1.89%  '
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = '0.592'
$ws.Range('E8').Value = '  +2.08%  '
$ws.Range('D9').Value = '3.343.30'
$ws.Range('E9').Value = '  +2.83%  '
$ws.Range('E10').Value = '  +6.00%  '
$ws.Range('D11').Value = '0.582'
$ws.Range('E11').Value = '  +2.86%  '
$ws.Range('D12').Value = '46.94'
$ws.Range('E12').Value = '  +4.46%  '
$ws.Range('D13').Value = '0.0000274'
$ws.Range('E13').Value = '  +3.02%  '
$ws.Range('D14').Value = '688.96'
$ws.Range('E14').Value = '  -0.13%  '
$ws.Range('D15').Value = '3.882.58'
$ws.Range('E15').Value = '  +2.76%  '
$ws.Range('D16').Value = '8.46'
$ws.Range('E16').Value = '  +2.68%  '
$ws.Range('D17').Value = '68.036.69'
$ws.Range('E17').Value = '  +2.13%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.345.83'
$ws.Range('E18').Value = '  +2.64%  '
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').Value = '0.118'
$ws.Range('E19').Value = '  -0.10%  '
$ws.Range('D20').Value = '17.44'
$ws.Range('E20').Value = '  +1.51%  '
$ws.Range('E21').Value = '  +4.44%  '
$ws.Range('D22').Value = '0.898'
$ws.Range('E22').Value = '  +2.08%  '
$ws.Range('E23').Value = '  +5.47%  '
$ws.Range('D24').Value = '17.11'
$ws.Range('E24').Value = '  +2.11%  '
$ws.Range('D25').Value = '98.61'
$ws.Range('E25').Value = '  +1.49%  '
$ws.Range('D26').Value = '3.90'
$ws.Range('E26').Value = '  +1.67%  '
$ws.Range('E27').Value = '  +0.87%  '
$ws.Range('D28').Value = '9.55'
$ws.Range('E28').Value = '  +3.87%  '
$ws.Range('D29').Value = '33.06'
$ws.Range('E29').Value = '  +1.59%  '
$ws.Range('D30').Value = '8.60'
$ws.Range('E30').Value = '  +3.33%  '
$ws.Range('D31').Value = '7.12'
$ws.Range('E31').Value = '  +6.94%  '
$ws.Range('D32').Value = '576.01'
$ws.Range('E32').Value = '  +0.54%  '
$ws.Range('D33').Value = '11.03'
$ws.Range('E33').Value = '  +3.19%  '
$ws.Range('E34').Value = '  +3.33%  '
$ws.Range('D35').Value = '3.728.82'
$ws.Range('E35').Value = '  -3.01%  '
$ws.Range('D36').Value = '57.40'
$ws.Range('E36').Value = '  +4.08%  '
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('D38').Value = '3.37'
$ws.Range('E38').Value = '  +3.95%  '
$ws.Range('D39').Value = '34.56'
$ws.Range('E39').Value = '  +9.85%  '
$ws.Range('E40').Value = '  +2.66%  '
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').Value = '2.67'
$ws.Range('E41').Value = '  +3.49%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '3.20'
$ws.Range('E42').Value = '  +6.09%  '
$ws.Range('D43').Value = '0.0₃0681'
$ws.Range('E43').Value = '  +3.21%  '
$ws.Range('D44').Value = '3.34'
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('D45').Value = '0.337'
$ws.Range('E45').Value = '  +4.13%  '
$ws.Range('D46').Value = '0.0407'
$ws.Range('E46').Value = '  +1.32%  '
$ws.Range('D47').Value = '2.67'
$ws.Range('E47').Value = '  +7.14%  '
$ws.Range('E48').Value = '  +1.95%  '
$ws.Range('E49').Value = '  -0.53%  '
$ws.Range('D50').Value = '1.33'
$ws.Range('E50').Value = '  -2.19%  '
$ws.Range('D51').Value = '129.42'
$ws.Range('E51').Value = '  +0.28%  '
